$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TestCases")
$ws1.Activate()

# The "Results" column (J) on rows 2-4 used to hold a "Pass" value for each
# test case. Those results (and the now-unused "Fail" shared string) are no
# longer wanted, so clear the cells out.
$ws1.Range("J2").Value = $null
$ws1.Range("J3").Value = $null
$ws1.Range("J4").Value = $null

# Move the selection / scroll position: was H4 (scrolled so row 2 is the
# top-left visible row) -> now just C1 with no special scroll offset.
$ws1.Range("C1").Select()
